$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The alcohol measurement data on Sheet1 had an extra data column (M)
# that needs to be dropped: deleting it discards its values and shifts
# everything to the right (the old column N) one step left, so the old
# column N becomes the new, last column M.
$ws.Columns("M").Delete() | Out-Null

# Park the selection on the column that is now the rightmost data
# column, mirroring where the deleted column used to live.
$ws.Range("M1").Select() | Out-Null
